# Apply the Italian-localization edits to the ISO-27005 5x5 risk matrix workbook.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # library_content
$ws2 = $wb.Worksheets.Item(2)   # spec

# ---------------------------------------------------------------------------
# Sheet "library_content" (sheet1)
# ---------------------------------------------------------------------------

# B2 (library_version) bumped from 3 to 4
$ws1.Range("B2").Value = 4

# New rows 19-22: Italian library/risk-matrix name+description pairs
$ws1.Range("A19").Value = "library_name[it]"
$ws1.Range("A19").Font.Name = "Aptos Narrow"

$ws1.Range("B19").Value = "Matrice 5x5 ISO-27005"
$ws1.Range("B19").HorizontalAlignment = -4131
$ws1.Range("B19").Font.Name = "Aptos Narrow"

$ws1.Range("A20").Value = "library_description[it]"
$ws1.Range("A20").Font.Name = "Aptos Narrow"

$ws1.Range("B20").Value = "Matrice 5x5 descritta nella norma ISO-27005 allegato A"
$ws1.Range("B20").HorizontalAlignment = -4131
$ws1.Range("B20").Font.Name = "Aptos Narrow"

$ws1.Range("A21").Value = "risk_matrix_name[it]"
$ws1.Range("A21").Font.Name = "Aptos Narrow"

$ws1.Range("B21").Value = "Matrice 5x5 ISO-27005"
$ws1.Range("B21").HorizontalAlignment = -4131
$ws1.Range("B21").Font.Name = "Aptos Narrow"

$ws1.Range("A22").Value = "risk_matrix_description[it]"
$ws1.Range("A22").Font.Name = "Aptos Narrow"

$ws1.Range("B22").Value = "Matrice 5x5 descritta nella norma ISO-27005 allegato A"
$ws1.Range("B22").HorizontalAlignment = -4131
$ws1.Range("B22").Font.Name = "Aptos Narrow"

# ---------------------------------------------------------------------------
# Sheet "spec" (sheet2): add name[it] / description[it] columns (N, O)
# ---------------------------------------------------------------------------

$ws2.Range("N1").Value = "name[it]"
$ws2.Range("O1").Value = "description[it]"

$ws2.Range("N2").Value = "5 - quasi certo"
$ws2.Range("N2").Font.Name = "Aptos Narrow"
$ws2.Range("O2").Value = "La verosimiglianza dello scenario di rischio è molto alta"
$ws2.Range("O2").Font.Name = "Aptos Narrow"

$ws2.Range("N3").Value = "4 - molto probabile"
$ws2.Range("N3").Font.Name = "Aptos Narrow"
$ws2.Range("O3").Value = "La verosimiglianza dello scenario di rischio è alta"
$ws2.Range("O3").Font.Name = "Aptos Narrow"

$ws2.Range("N4").Value = "3 - probabile"
$ws2.Range("N4").Font.Name = "Aptos Narrow"
$ws2.Range("O4").Value = "La verosimiglianza dello scenario di rischio è significativa"
$ws2.Range("O4").Font.Name = "Aptos Narrow"

$ws2.Range("N5").Value = "2 - piuttosto improbabile"
$ws2.Range("N5").Font.Name = "Aptos Narrow"
$ws2.Range("O5").Value = "La verosimiglianza dello scenario di rischio è bassa"
$ws2.Range("O5").Font.Name = "Aptos Narrow"

$ws2.Range("N6").Value = "1 - improbabile"
$ws2.Range("N6").Font.Name = "Aptos Narrow"
$ws2.Range("O6").Value = "La verosimiglianza dello scenario di rischio è molto bassa"
$ws2.Range("O6").Font.Name = "Aptos Narrow"

$ws2.Range("N7").Value = "1 - minore"
$ws2.Range("N7").Font.Name = "Aptos Narrow"
$ws2.Range("O7").Value = "Conseguenze trascurabili per l'organizzazione"
$ws2.Range("O7").Font.Name = "Aptos Narrow"

$ws2.Range("N8").Value = "2 - significativo"
$ws2.Range("N8").Font.Name = "Aptos Narrow"
$ws2.Range("O8").Value = "Conseguenze significative ma limitate per l'organizzazione"
$ws2.Range("O8").Font.Name = "Aptos Narrow"

$ws2.Range("N9").Value = "3 - serio"
$ws2.Range("N9").Font.Name = "Aptos Narrow"
$ws2.Range("O9").Value = "Conseguenze sostanziali per l'organizzazione"
$ws2.Range("O9").Font.Name = "Aptos Narrow"

$ws2.Range("N10").Value = "4 - critico"
$ws2.Range("N10").Font.Name = "Aptos Narrow"
$ws2.Range("O10").Value = "Conseguenze disastrose per l'organizzazione"
$ws2.Range("O10").Font.Name = "Aptos Narrow"

$ws2.Range("N11").Value = "5 - catastrofico"
$ws2.Range("N11").Font.Name = "Aptos Narrow"
$ws2.Range("O11").Value = "Conseguenze settoriali o regolamentari oltre l'organizzazione"
$ws2.Range("O11").Font.Name = "Aptos Narrow"

$ws2.Range("N12").Value = "1 - molto basso"
$ws2.Range("N12").Font.Name = "Aptos Narrow"
$ws2.Range("O12").Value = "molto basso - rischio accettabile"
$ws2.Range("O12").Font.Name = "Aptos Narrow"

$ws2.Range("N13").Value = "2 - basso"
$ws2.Range("N13").Font.Name = "Aptos Narrow"
$ws2.Range("O13").Value = "basso - rischio accettabile"
$ws2.Range("O13").Font.Name = "Aptos Narrow"

$ws2.Range("N14").Value = "3 - medio"
$ws2.Range("N14").Font.Name = "Aptos Narrow"
$ws2.Range("O14").Value = "medio - rischio tollerabile"
$ws2.Range("O14").Font.Name = "Aptos Narrow"

$ws2.Range("N15").Value = "4 - alto"
$ws2.Range("N15").Font.Name = "Aptos Narrow"
$ws2.Range("O15").Value = "alto - rischio inaccettabile"
$ws2.Range("O15").Font.Name = "Aptos Narrow"

$ws2.Range("N16").Value = "5 - molto alto"
$ws2.Range("N16").Font.Name = "Aptos Narrow"
$ws2.Range("O16").Value = "molto alto - rischio inaccettabile"
$ws2.Range("O16").Font.Name = "Aptos Narrow"

# Column widths for the newly populated columns (best-fit approximations)
$ws2.Columns.Item(13).ColumnWidth = 58.33
$ws2.Columns.Item(14).ColumnWidth = 20.17
$ws2.Columns.Item(15).ColumnWidth = 50.33

# ---------------------------------------------------------------------------
# View state: selections + zoom (applied last so the active tab ends on sheet1)
# ---------------------------------------------------------------------------

$ws2.Select()
$ws2.Range("Q3").Select()
$excel.ActiveWindow.Zoom = 125

$ws1.Select()
$ws1.Range("B2").Select()

Write-Host "Edit complete"
